$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new data values for C10 (bennett) and C15 (frederick)
$ws.Range("C10").Value = 2010
$ws.Range("C15").Value = 470

# Update the active cell selection to C17
$ws.Range("C17").Select()
